$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.077022314071655
$ws.Range("B1").Value = 6.279294967651367
$ws.Range("C1").Value = 6.324397563934326
$ws.Range("D1").Value = 6.721923828125
$ws.Range("E1").Value = 5.208512306213379
